# Update countries & provincias Spain
# Refresh COVID case figures and re-sort a few countries that moved in the
# total-cases ranking (Tailandia above Egipto; Kazajistan above Estonia).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Reorder countries in column A (Tailandia/Egipto swap; Kazajistan moved up) ---
$ws.Range("A53").Value = "Tailandia"
$ws.Range("A54").Value = "Egipto"

$ws.Range("A67").Value = "Kazajistan"
$ws.Range("A68").Value = "Estonia"
$ws.Range("A69").Value = "Irak"
$ws.Range("A70").Value = "Nueva Zelanda"

# --- Update case numbers (columns B-H) for affected rows ---
$ws.Range("D8").Value = 81800
$ws.Range("E8").Value = 51846

$ws.Range("B14").Value = 30891
$ws.Range("C14").Value = 208
$ws.Range("E14").Value = 14913
$ws.Range("G14").Value = 5
$ws.Range("H14").Value = 1952

$ws.Range("B32").Value = 7025
$ws.Range("C32").Value = 106
$ws.Range("D32").Value = 1765
$ws.Range("E32").Value = 5125
$ws.Range("G32").Value = 7
$ws.Range("H32").Value = 135

$ws.Range("B53").Value = 2700
$ws.Range("C53").Value = 28
$ws.Range("D53").Value = 1689
$ws.Range("E53").Value = 964
$ws.Range("F53").Value = 61
$ws.Range("G53").Value = 1
$ws.Range("H53").Value = 47

$ws.Range("B54").Value = 2673
$ws.Range("D54").Value = 596
$ws.Range("E54").Value = 1881
$ws.Range("F54").Value = 0
$ws.Range("H54").Value = 196

$ws.Range("B67").Value = 1470
$ws.Range("C67").Value = 68
$ws.Range("D67").Value = 277
$ws.Range("E67").Value = 1176
$ws.Range("F67").Value = 22
$ws.Range("H67").Value = 17

$ws.Range("D68").Value = 133
$ws.Range("E68").Value = 1265
$ws.Range("F68").Value = 10
$ws.Range("H68").Value = 36

$ws.Range("B69").Value = 1434
$ws.Range("C69").Value = 0
$ws.Range("D69").Value = 856
$ws.Range("E69").Value = 498
$ws.Range("F69").Value = 0
$ws.Range("G69").Value = 0
$ws.Range("H69").Value = 80

$ws.Range("B70").Value = 1409
$ws.Range("C70").Value = 8
$ws.Range("D70").Value = 816
$ws.Range("E70").Value = 582
$ws.Range("F70").Value = 2
$ws.Range("G70").Value = 2
$ws.Range("H70").Value = 11

$ws.Range("E78").Value = 838
$ws.Range("G78").Value = 1
$ws.Range("H78").Value = 5

$ws.Range("B99").Value = 489
$ws.Range("C99").Value = 23
$ws.Range("D99").Value = 114

$ws.Range("D115").Value = 194
$ws.Range("E115").Value = 74
